$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column header from "name" to "display_name" (part of variable name unification)
$ws.Range("A1").Value = "display_name"

# Resize column A to fit the new, longer header text
$ws.Columns.Item(1).ColumnWidth = 11.8

# Update the active cell selection to match the edited worksheet state
[void]$ws.Range("D11").Select()
